# Update the cryptos price list with freshly scraped values.
# D-column values that parse as plain decimal numbers (e.g. "287.25") are
# given a leading apostrophe so Excel stores them as text (matching the
# original inlineStr cells) instead of silently converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.075.39"
$ws.Range("E2").Value = "  -1.71%  "
$ws.Range("D3").Value = "1.550.15"
$ws.Range("E3").Value = "  -1.27%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("E5").Value = "  -0.03%  "
$ws.Range("D6").Value = "'287.25"
$ws.Range("E6").Value = "  -0.37%  "
$ws.Range("D7").Value = "'0.3829"
$ws.Range("E7").Value = "  +2.85%  "
$ws.Range("D8").Value = "'0.3276"
$ws.Range("E8").Value = "  -1.24%  "
$ws.Range("D9").Value = "'43.58"
$ws.Range("E9").Value = "  -9.73%  "
$ws.Range("D10").Value = "'1.121"
$ws.Range("E10").Value = "  -1.14%  "
$ws.Range("D11").Value = "'0.07347"
$ws.Range("E11").Value = "  -1.79%  "
$ws.Range("E12").Value = "  +0.01%  "
$ws.Range("D13").Value = "'19.96"
$ws.Range("E13").Value = "  -3.53%  "
$ws.Range("D14").Value = "'5.783"
$ws.Range("E14").Value = "  -2.59%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "'6.743"
$ws.Range("E15").Value = "  -2.38%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "1.566.47"
$ws.Range("E16").Value = "  +0.01%  "
$ws.Range("D17").Value = "'0.00001081"
$ws.Range("E17").Value = "  -3.21%  "
$ws.Range("D18").Value = "'0.06616"
$ws.Range("E18").Value = "  -1.89%  "
$ws.Range("D19").Value = "'85.50"
$ws.Range("E19").Value = "  -2.71%  "
$ws.Range("E20").Value = "  -0.12%  "
$ws.Range("D21").Value = "'6.365"
$ws.Range("E21").Value = "  +0.15%  "
$ws.Range("D22").Value = "'15.98"
$ws.Range("E22").Value = "  -3.08%  "
$ws.Range("D23").Value = "'11.64"
$ws.Range("E23").Value = "  -3.58%  "
$ws.Range("D24").Value = "22.069.06"
$ws.Range("E24").Value = "  -1.73%  "
$ws.Range("D25").Value = "'2.294"
$ws.Range("E25").Value = "  -3.98%  "
$ws.Range("D26").Value = "'2.483"
$ws.Range("E26").Value = "  -3.23%  "
$ws.Range("D27").Value = "'150.01"
$ws.Range("E27").Value = "  -2.01%  "
$ws.Range("D28").Value = "'19.05"
$ws.Range("E28").Value = "  -3.43%  "
$ws.Range("D29").Value = "'4.930"
$ws.Range("E29").Value = "  -1.87%  "
$ws.Range("D30").Value = "1.750.46"
$ws.Range("E30").Value = "  +0.40%  "
$ws.Range("D31").Value = "'121.04"
$ws.Range("E31").Value = "  -2.63%  "
$ws.Range("D32").Value = "'1.072"
$ws.Range("E32").Value = "  +1.76%  "
$ws.Range("D33").Value = "'5.856"
$ws.Range("E33").Value = "  -4.49%  "
$ws.Range("D34").Value = "'1.898"
$ws.Range("E34").Value = "  -5.54%  "
$ws.Range("D35").Value = "'0.08223"
$ws.Range("E35").Value = "  -1.22%  "
$ws.Range("D36").Value = "'9.187"
$ws.Range("E36").Value = "  -5.95%  "
$ws.Range("D37").Value = "'0.06275"
$ws.Range("E37").Value = "  -1.91%  "
$ws.Range("D38").Value = "'0.02306"
$ws.Range("E38").Value = "  -6.41%  "
$ws.Range("D39").Value = "'5.254"
$ws.Range("E39").Value = "  -2.35%  "
$ws.Range("D40").Value = "'0.2148"
$ws.Range("E40").Value = "  -5.54%  "
$ws.Range("D41").Value = "'1.231"
$ws.Range("E41").Value = "  -4.20%  "
$ws.Range("D42").Value = "'10.94"
$ws.Range("E42").Value = "  -3.19%  "
$ws.Range("E43").Value = "  -0.05%  "
$ws.Range("D44").Value = "'0.5995"
$ws.Range("E44").Value = "  -4.97%  "
$ws.Range("D45").Value = "'13.73"
$ws.Range("E45").Value = "  -0.86%  "
$ws.Range("D46").Value = "'3.729"
$ws.Range("E46").Value = "  -1.27%  "
$ws.Range("D47").Value = "'0.5783"
$ws.Range("E47").Value = "  -6.06%  "
$ws.Range("D48").Value = "'1.965"
$ws.Range("E48").Value = "  -4.29%  "
$ws.Range("D49").Value = "'121.62"
$ws.Range("E49").Value = "  -3.21%  "
$ws.Range("D50").Value = "'1.172"
$ws.Range("E50").Value = "  -3.18%  "
$ws.Range("D51").Value = "'0.07011"
$ws.Range("E51").Value = "  -2.83%  "
